# Update crypto price/volume data per latest scrape (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.875.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.17%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.158.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.10%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.13%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.42%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.156.83"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.15%  "

$ws.Range("E9").Value = "  +4.35%  "

$ws.Range("E10").Value = "  +5.98%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.19"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.52%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.504"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.90%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000256"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +12.35%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.63%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.679.20"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.41%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.997.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.31%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.89%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.160.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.11%  "

$ws.Range("E19").Value = "  +0.50%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "513.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.13%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.90%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.736"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.52%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.85"
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.89%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.47%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.80%  "

$ws.Range("E29").Value = "  +6.15%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.57%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.11%  "

$ws.Range("E32").Value = "  +3.92%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.65"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.88%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.91%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.63"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "487.18"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.92%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0424"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.78%  "

$ws.Range("E40").Value = "  +1.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.120.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.52%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.64%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.121"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.08%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.294"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.47%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +15.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.66%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0583"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +12.70%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").Style = "Normal"

$ws.Range("E49").Value = "  +3.29%  "

$ws.Range("E50").Value = "  +10.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "119.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.93%  "

